# edit.ps1 - applies the 2022-Q3 sheet insertion + 总计 sheet updates
$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, [string]$val) {
    # Force a cell to store $val verbatim as text (preserve leading/trailing zeros),
    # then reset number formatting back to the default style so no stray style
    # index (quotePrefix / custom numFmt) is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---- 1. "总计" (summary) sheet: shift existing rows down, add 2022-Q3 totals ----
$summary = $wb.Worksheets.Item(1)

$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 24
$summary.Cells.Item(2,4).Value = 9.039999999999999

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q2"
$summary.Cells.Item(3,3).Value = 17
$summary.Cells.Item(3,4).Value = 2.8

# Row 4 ("2022-Q1") is a brand-new row in the grid; copy the style used by the
# existing A-column cells (bold/centered/bordered) onto its A cell before writing it.
$summary.Cells.Item(2,1).Copy()
$summary.Cells.Item(4,1).PasteSpecial(-4122)
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2022-Q1"
$summary.Cells.Item(4,3).Value = 7
$summary.Cells.Item(4,4).Value = 1.47

# ---- 2. Insert a brand-new "2022-Q3" sheet right after "总计" ----
$new = $wb.Worksheets.Add($null, $summary)
$new.Name = "2022-Q3"

# Match the page-margin defaults used by the rest of the workbook's sheets
# (0.75in/1in/0.5in) instead of the engine's own Excel-style defaults.
$new.PageSetup.LeftMargin = 54
$new.PageSetup.RightMargin = 54
$new.PageSetup.TopMargin = 72
$new.PageSetup.BottomMargin = 72
$new.PageSetup.HeaderMargin = 36
$new.PageSetup.FooterMargin = 36

# Header row (basic labels) - values first, then copy the header style
# (bold font + border, centered) from the equivalent header on the "2022-Q2" sheet.
$oldQ2 = $wb.Worksheets.Item(3)
$new.Cells.Item(1,2).Value = "基金代码"
$new.Cells.Item(1,3).Value = "基金名称"
$new.Cells.Item(1,4).Value = "基金规模"
$new.Cells.Item(1,5).Value = "股票总仓位"
$new.Cells.Item(1,6).Value = "仓位占比"
$new.Cells.Item(1,7).Value = "持有市值(亿元)"
$new.Cells.Item(1,8).Value = "仓位排名"
$oldQ2.Range("B1:H1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)

Set-TextCell $new.Cells.Item(2,2) "002251"
$new.Cells.Item(2,3).Value = "华夏军工安全灵活配置混合A"
Set-TextCell $new.Cells.Item(2,4) "42.95"
Set-TextCell $new.Cells.Item(2,5) "94.59"
Set-TextCell $new.Cells.Item(2,6) "4.61"
Set-TextCell $new.Cells.Item(2,7) "1.9800"
$new.Cells.Item(2,1).Value = 0
$new.Cells.Item(2,8).Value = 9
Set-TextCell $new.Cells.Item(3,2) "005774"
$new.Cells.Item(3,3).Value = "华夏产业升级混合A"
Set-TextCell $new.Cells.Item(3,4) "24.29"
Set-TextCell $new.Cells.Item(3,5) "93.85"
Set-TextCell $new.Cells.Item(3,6) "7.62"
Set-TextCell $new.Cells.Item(3,7) "1.8509"
$new.Cells.Item(3,1).Value = 1
$new.Cells.Item(3,8).Value = 3
Set-TextCell $new.Cells.Item(4,2) "012390"
$new.Cells.Item(4,3).Value = "中欧产业前瞻混合A"
Set-TextCell $new.Cells.Item(4,4) "16.91"
Set-TextCell $new.Cells.Item(4,5) "91.97"
Set-TextCell $new.Cells.Item(4,6) "4.69"
Set-TextCell $new.Cells.Item(4,7) "0.7931"
$new.Cells.Item(4,1).Value = 2
$new.Cells.Item(4,8).Value = 8
Set-TextCell $new.Cells.Item(5,2) "010481"
$new.Cells.Item(5,3).Value = "汇添富高质量成长精选2年持有期混合"
Set-TextCell $new.Cells.Item(5,4) "31.12"
Set-TextCell $new.Cells.Item(5,5) "76.14"
Set-TextCell $new.Cells.Item(5,6) "2.29"
Set-TextCell $new.Cells.Item(5,7) "0.7126"
$new.Cells.Item(5,1).Value = 3
$new.Cells.Item(5,8).Value = 9
Set-TextCell $new.Cells.Item(6,2) "015059"
$new.Cells.Item(6,3).Value = "华夏产业升级混合C"
Set-TextCell $new.Cells.Item(6,4) "8.92"
Set-TextCell $new.Cells.Item(6,5) "93.85"
Set-TextCell $new.Cells.Item(6,6) "7.62"
Set-TextCell $new.Cells.Item(6,7) "0.6797"
$new.Cells.Item(6,1).Value = 4
$new.Cells.Item(6,8).Value = 3
Set-TextCell $new.Cells.Item(7,2) "013566"
$new.Cells.Item(7,3).Value = "华夏军工安全灵活配置混合C"
Set-TextCell $new.Cells.Item(7,4) "11.35"
Set-TextCell $new.Cells.Item(7,5) "94.59"
Set-TextCell $new.Cells.Item(7,6) "4.61"
Set-TextCell $new.Cells.Item(7,7) "0.5232"
$new.Cells.Item(7,1).Value = 5
$new.Cells.Item(7,8).Value = 9
Set-TextCell $new.Cells.Item(8,2) "012155"
$new.Cells.Item(8,3).Value = "汇添富成长先锋六个月持有期混合A"
Set-TextCell $new.Cells.Item(8,4) "15.78"
Set-TextCell $new.Cells.Item(8,5) "85.82"
Set-TextCell $new.Cells.Item(8,6) "3.21"
Set-TextCell $new.Cells.Item(8,7) "0.5065"
$new.Cells.Item(8,1).Value = 6
$new.Cells.Item(8,8).Value = 8
Set-TextCell $new.Cells.Item(9,2) "012557"
$new.Cells.Item(9,3).Value = "中欧景气前瞻一年持有期混合型证券投资基金A"
Set-TextCell $new.Cells.Item(9,4) "8.43"
Set-TextCell $new.Cells.Item(9,5) "92.09"
Set-TextCell $new.Cells.Item(9,6) "4.50"
Set-TextCell $new.Cells.Item(9,7) "0.3794"
$new.Cells.Item(9,1).Value = 7
$new.Cells.Item(9,8).Value = 8
Set-TextCell $new.Cells.Item(10,2) "010599"
$new.Cells.Item(10,3).Value = "汇添富高质量成长30一年持有期混合A"
Set-TextCell $new.Cells.Item(10,4) "15.29"
Set-TextCell $new.Cells.Item(10,5) "85.27"
Set-TextCell $new.Cells.Item(10,6) "2.25"
Set-TextCell $new.Cells.Item(10,7) "0.3440"
$new.Cells.Item(10,1).Value = 8
$new.Cells.Item(10,8).Value = 10
Set-TextCell $new.Cells.Item(11,2) "005358"
$new.Cells.Item(11,3).Value = "东方阿尔法精选灵活配置混合A"
Set-TextCell $new.Cells.Item(11,4) "3.47"
Set-TextCell $new.Cells.Item(11,5) "93.91"
Set-TextCell $new.Cells.Item(11,6) "9.81"
Set-TextCell $new.Cells.Item(11,7) "0.3404"
$new.Cells.Item(11,1).Value = 9
$new.Cells.Item(11,8).Value = 1
Set-TextCell $new.Cells.Item(12,2) "010615"
$new.Cells.Item(12,3).Value = "国金自主创新混合A"
Set-TextCell $new.Cells.Item(12,4) "3.42"
Set-TextCell $new.Cells.Item(12,5) "82.69"
Set-TextCell $new.Cells.Item(12,6) "5.58"
Set-TextCell $new.Cells.Item(12,7) "0.1908"
$new.Cells.Item(12,1).Value = 10
$new.Cells.Item(12,8).Value = 8
Set-TextCell $new.Cells.Item(13,2) "010616"
$new.Cells.Item(13,3).Value = "国金自主创新混合C"
Set-TextCell $new.Cells.Item(13,4) "3.42"
Set-TextCell $new.Cells.Item(13,5) "82.69"
Set-TextCell $new.Cells.Item(13,6) "5.58"
Set-TextCell $new.Cells.Item(13,7) "0.1908"
$new.Cells.Item(13,1).Value = 11
$new.Cells.Item(13,8).Value = 8
Set-TextCell $new.Cells.Item(14,2) "001173"
$new.Cells.Item(14,3).Value = "中欧瑾和灵活配置混合 - A"
Set-TextCell $new.Cells.Item(14,4) "2.26"
Set-TextCell $new.Cells.Item(14,5) "92.00"
Set-TextCell $new.Cells.Item(14,6) "5.47"
Set-TextCell $new.Cells.Item(14,7) "0.1236"
$new.Cells.Item(14,1).Value = 12
$new.Cells.Item(14,8).Value = 5
Set-TextCell $new.Cells.Item(15,2) "001760"
$new.Cells.Item(15,3).Value = "嘉实创新成长灵活配置混合"
Set-TextCell $new.Cells.Item(15,4) "1.03"
Set-TextCell $new.Cells.Item(15,5) "93.49"
Set-TextCell $new.Cells.Item(15,6) "9.27"
Set-TextCell $new.Cells.Item(15,7) "0.0955"
$new.Cells.Item(15,1).Value = 13
$new.Cells.Item(15,8).Value = 2
Set-TextCell $new.Cells.Item(16,2) "014818"
$new.Cells.Item(16,3).Value = "国金新兴价值混合A"
Set-TextCell $new.Cells.Item(16,4) "1.25"
Set-TextCell $new.Cells.Item(16,5) "81.81"
Set-TextCell $new.Cells.Item(16,6) "5.32"
Set-TextCell $new.Cells.Item(16,7) "0.0665"
$new.Cells.Item(16,1).Value = 14
$new.Cells.Item(16,8).Value = 9
Set-TextCell $new.Cells.Item(17,2) "006803"
$new.Cells.Item(17,3).Value = "嘉实互通精选股票"
Set-TextCell $new.Cells.Item(17,4) "1.25"
Set-TextCell $new.Cells.Item(17,5) "91.91"
Set-TextCell $new.Cells.Item(17,6) "4.98"
Set-TextCell $new.Cells.Item(17,7) "0.0622"
$new.Cells.Item(17,1).Value = 15
$new.Cells.Item(17,8).Value = 7
Set-TextCell $new.Cells.Item(18,2) "005359"
$new.Cells.Item(18,3).Value = "东方阿尔法精选灵活配置混合C"
Set-TextCell $new.Cells.Item(18,4) "0.58"
Set-TextCell $new.Cells.Item(18,5) "93.91"
Set-TextCell $new.Cells.Item(18,6) "9.81"
Set-TextCell $new.Cells.Item(18,7) "0.0569"
$new.Cells.Item(18,1).Value = 16
$new.Cells.Item(18,8).Value = 1
Set-TextCell $new.Cells.Item(19,2) "012558"
$new.Cells.Item(19,3).Value = "中欧景气前瞻一年持有期混合型证券投资基金C"
Set-TextCell $new.Cells.Item(19,4) "0.68"
Set-TextCell $new.Cells.Item(19,5) "92.09"
Set-TextCell $new.Cells.Item(19,6) "4.50"
Set-TextCell $new.Cells.Item(19,7) "0.0306"
$new.Cells.Item(19,1).Value = 17
$new.Cells.Item(19,8).Value = 8
Set-TextCell $new.Cells.Item(20,2) "014819"
$new.Cells.Item(20,3).Value = "国金新兴价值混合C"
Set-TextCell $new.Cells.Item(20,4) "0.54"
Set-TextCell $new.Cells.Item(20,5) "81.81"
Set-TextCell $new.Cells.Item(20,6) "5.32"
Set-TextCell $new.Cells.Item(20,7) "0.0287"
$new.Cells.Item(20,1).Value = 18
$new.Cells.Item(20,8).Value = 9
Set-TextCell $new.Cells.Item(21,2) "012391"
$new.Cells.Item(21,3).Value = "中欧产业前瞻混合C"
Set-TextCell $new.Cells.Item(21,4) "0.52"
Set-TextCell $new.Cells.Item(21,5) "91.97"
Set-TextCell $new.Cells.Item(21,6) "4.69"
Set-TextCell $new.Cells.Item(21,7) "0.0244"
$new.Cells.Item(21,1).Value = 19
$new.Cells.Item(21,8).Value = 8
Set-TextCell $new.Cells.Item(22,2) "011259"
$new.Cells.Item(22,3).Value = "汇添富高质量成长30一年持有期混合C"
Set-TextCell $new.Cells.Item(22,4) "0.70"
Set-TextCell $new.Cells.Item(22,5) "85.27"
Set-TextCell $new.Cells.Item(22,6) "2.25"
Set-TextCell $new.Cells.Item(22,7) "0.0158"
$new.Cells.Item(22,1).Value = 20
$new.Cells.Item(22,8).Value = 10
Set-TextCell $new.Cells.Item(23,2) "012156"
$new.Cells.Item(23,3).Value = "汇添富成长先锋六个月持有期混合C"
Set-TextCell $new.Cells.Item(23,4) "0.44"
Set-TextCell $new.Cells.Item(23,5) "85.82"
Set-TextCell $new.Cells.Item(23,6) "3.21"
Set-TextCell $new.Cells.Item(23,7) "0.0141"
$new.Cells.Item(23,1).Value = 21
$new.Cells.Item(23,8).Value = 8
Set-TextCell $new.Cells.Item(24,2) "002567"
$new.Cells.Item(24,3).Value = "大成国家安全主题灵活配置混合"
Set-TextCell $new.Cells.Item(24,4) "0.35"
Set-TextCell $new.Cells.Item(24,5) "71.95"
Set-TextCell $new.Cells.Item(24,6) "4.01"
Set-TextCell $new.Cells.Item(24,7) "0.0140"
$new.Cells.Item(24,1).Value = 22
$new.Cells.Item(24,8).Value = 8
Set-TextCell $new.Cells.Item(25,2) "001174"
$new.Cells.Item(25,3).Value = "中欧瑾和灵活配置混合 - C"
Set-TextCell $new.Cells.Item(25,4) "0.23"
Set-TextCell $new.Cells.Item(25,5) "92.00"
Set-TextCell $new.Cells.Item(25,6) "5.47"
Set-TextCell $new.Cells.Item(25,7) "0.0126"
$new.Cells.Item(25,1).Value = 23
$new.Cells.Item(25,8).Value = 5

# Column-A numbering style (bold/centered/bordered), matching the other sheets,
# copied onto the whole A2:A25 data range in one paste.
$oldQ2.Cells.Item(2,1).Copy()
$new.Range("A2:A25").PasteSpecial(-4122)

# Re-enter the numeric index values (the PasteSpecial above only copied formats,
# but belt-and-braces re-assert the values in case a prior step left any blank).
$new.Cells.Item(2,1).Value = 0
$new.Cells.Item(3,1).Value = 1
$new.Cells.Item(4,1).Value = 2
$new.Cells.Item(5,1).Value = 3
$new.Cells.Item(6,1).Value = 4
$new.Cells.Item(7,1).Value = 5
$new.Cells.Item(8,1).Value = 6
$new.Cells.Item(9,1).Value = 7
$new.Cells.Item(10,1).Value = 8
$new.Cells.Item(11,1).Value = 9
$new.Cells.Item(12,1).Value = 10
$new.Cells.Item(13,1).Value = 11
$new.Cells.Item(14,1).Value = 12
$new.Cells.Item(15,1).Value = 13
$new.Cells.Item(16,1).Value = 14
$new.Cells.Item(17,1).Value = 15
$new.Cells.Item(18,1).Value = 16
$new.Cells.Item(19,1).Value = 17
$new.Cells.Item(20,1).Value = 18
$new.Cells.Item(21,1).Value = 19
$new.Cells.Item(22,1).Value = 20
$new.Cells.Item(23,1).Value = 21
$new.Cells.Item(24,1).Value = 22
$new.Cells.Item(25,1).Value = 23

# ---- 3. Restore the active tab to "2022-Q1" (now the last / 4th sheet), matching ----
# ---- the original workbook state where the last sheet was the selected tab.      ----
$wb.Worksheets.Item(4).Activate()

Write-Output "edit complete"
